$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-09 Thursday", "2024-05-10 Friday"),
    @("62×85=5270", "35×97=3395"),
    @("61×66=4026", "51×95=4845"),
    @("65×85=5525", "29×97=2813"),
    @("62×96=5952", "45×50=2250"),
    @("26×92=2392", "24×60=1440"),
    @("36×86=3096", "47×29=1363"),
    @("25×86=2150", "25×95=2375"),
    @("21×54=1134", "45×66=2970"),
    @("84×64=5376", "78×41=3198"),
    @("63×21=1323", "21×11=231"),
    @("46×66=3036", "98×99=9702"),
    @("93×93=8649", "16×75=1200"),
    @("33×53=1749", "94×67=6298"),
    @("65×26=1690", "33×26=858"),
    @("93×71=6603", "67×76=5092"),
    @("68×53=3604", "29×16=464"),
    @("92×54=4968", "55×82=4510"),
    @("65×86=5590", "81×22=1782"),
    @("32×20=640", "55×34=1870"),
    @("47×27=1269", "11×44=484"),
    @("76×20=1520", "99×73=7227"),
    @("43×47=2021", "18×27=486"),
    @("22×55=1210", "14×12=168"),
    @("78×54=4212", "51×20=1020"),
    @("82×89=7298", "79×72=5688")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
